$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price cells that look numeric, so Excel
# keeps them as text (matching the original inline-string cell type)
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D11", "D14", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D28", "D33", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.447.65'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '3.095.75'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '582.60'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").Value = '144.81'
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.085.80'
$ws.Range("E8").Value = '  -0.78%  '
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("E10").Value = '  +6.31%  '
$ws.Range("D11").Value = '5.59'
$ws.Range("E11").Value = '  -3.58%  '
$ws.Range("E12").Value = '  -2.61%  '
$ws.Range("E13").Value = '  -1.91%  '
$ws.Range("D14").Value = '37.11'
$ws.Range("E14").Value = '  +4.14%  '
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("D16").Value = '3.609.99'
$ws.Range("E16").Value = '  -0.71%  '
$ws.Range("D17").Value = '63.284.67'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '7.06'
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("D19").Value = '3.095.54'
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("D20").Value = '459.38'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").Value = '14.18'
$ws.Range("E21").Value = '  +0.85%  '
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").Value = '7.39'
$ws.Range("E23").Value = '  -1.83%  '
$ws.Range("D24").Value = '81.15'
$ws.Range("E24").Value = '  -1.13%  '
$ws.Range("D25").Value = '12.90'
$ws.Range("E25").Value = '  -3.11%  '
$ws.Range("E26").Value = '  -2.09%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '9.01'
$ws.Range("E28").Value = '  +8.77%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  -0.81%  '
$ws.Range("E31").Value = '  -2.37%  '
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("D33").Value = '0.111'
$ws.Range("E33").Value = '  +0.73%  '
$ws.Range("E34").Value = '  -1.65%  '
$ws.Range("E35").Value = '  -3.16%  '
$ws.Range("E36").Value = '  +2.14%  '
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("E38").Value = '  -5.23%  '
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("D40").Value = '50.23'
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("D41").Value = '431.63'
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").Value = '8.67'
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("D43").Value = '2.874.06'
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("D44").Value = '0.0365'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("D45").Value = '0.269'
$ws.Range("E45").Value = '  -3.63%  '
$ws.Range("E46").Value = '  -3.76%  '
$ws.Range("D47").Value = '35.82'
$ws.Range("E47").Value = '  +1.83%  '
$ws.Range("D48").Value = '124.64'
$ws.Range("E48").Value = '  +0.70%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  -1.42%  '
$ws.Range("D51").Value = '23.94'
$ws.Range("E51").Value = '  -2.45%  '
